# Add no-format options - extend README
#
# This script reproduces, via the Excel COM object model, the edit that:
#  - moves the old K1 ("foo") content out to N1
#  - introduces three new formatted cells: K1 (date DD/MM/YY), L1 (boolean
#    rendered as TRUE/FALSE), M1 (percentage 0.00%)
#  - bumps G1's value
#  - re-applies B1's (General) number format so it gets its own style entry
#  - updates the active selection / top-left cell of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the trailing "foo" label from K1 out to the new last column N1 ---
$ws.Range("N1").Value = "foo"

# --- B1 keeps its value, but gets an explicitly (re)applied General format ---
$ws.Range("B1").NumberFormat = "General"

# --- G1's value changes ---
$ws.Range("G1").Value = 200.666661562376

# --- K1 becomes a date, formatted DD/MM/YY ---
$ws.Range("K1").Value = 43911
$ws.Range("K1").NumberFormat = "DD/MM/YY"

# --- L1 becomes a boolean, formatted as TRUE/FALSE ---
$ws.Range("L1").Value = $true
$ws.Range("L1").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# --- M1 becomes a percentage value ---
$ws.Range("M1").Value = 1.23
$ws.Range("M1").NumberFormat = "0.00%"

# --- Update the view: scroll so column B is the left-most visible column,
#     and select M2 as the active cell ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("M2").Select()
